$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert 3 new rows before the last (bottom-bordered) data row (row 33)
# so the table grows from 18 data rows (16-33) to 21 data rows (16-36).
$ws.Rows("33:35").Insert()

# Step 2: Copy formatting from row 32 (a normal interior row) onto the new rows 33-35
# so they get the same cell styles as the rest of the table instead of Excel
# auto-generating new blended styles.
$ws.Range("B32:J32").Copy()
$ws.Range("B33:J35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Step 3: Clear all the old data values (formatting is preserved) so that the
# shared-string table is rebuilt cleanly when the workbook is saved.
$ws.Range("B16:J36").ClearContents()

# Step 4: Write the new data set (same 4 workers, now with a fuller set of periods)
$ws.Range("B16").Value2 = "CC"
$ws.Range("C16").Value2 = "73201153"
$ws.Range("D16").Value2 = "NELSON ALFONSO MENDOZA DIAZ"
$ws.Range("E16").Value2 = "1910"
$ws.Range("F16").Value2 = 33125
$ws.Range("G16").Value2 = 828116
$ws.Range("B17").Value2 = "CC"
$ws.Range("C17").Value2 = "73201153"
$ws.Range("D17").Value2 = "NELSON ALFONSO MENDOZA DIAZ"
$ws.Range("E17").Value2 = "1909"
$ws.Range("F17").Value2 = 33125
$ws.Range("G17").Value2 = 828116
$ws.Range("B18").Value2 = "CC"
$ws.Range("C18").Value2 = "8373999"
$ws.Range("D18").Value2 = "ELVER JOSE SEHUANES BULLOSO"
$ws.Range("E18").Value2 = "2003"
$ws.Range("F18").Value2 = 32021
$ws.Range("G18").Value2 = 828116
$ws.Range("B19").Value2 = "CC"
$ws.Range("C19").Value2 = "8373999"
$ws.Range("D19").Value2 = "ELVER JOSE SEHUANES BULLOSO"
$ws.Range("E19").Value2 = "2002"
$ws.Range("F19").Value2 = 33125
$ws.Range("G19").Value2 = 828116
$ws.Range("B20").Value2 = "CC"
$ws.Range("C20").Value2 = "8373999"
$ws.Range("D20").Value2 = "ELVER JOSE SEHUANES BULLOSO"
$ws.Range("E20").Value2 = "2001"
$ws.Range("F20").Value2 = 33125
$ws.Range("G20").Value2 = 828116
$ws.Range("B21").Value2 = "CC"
$ws.Range("C21").Value2 = "8373999"
$ws.Range("D21").Value2 = "ELVER JOSE SEHUANES BULLOSO"
$ws.Range("E21").Value2 = "1912"
$ws.Range("F21").Value2 = 33125
$ws.Range("G21").Value2 = 828116
$ws.Range("B22").Value2 = "CC"
$ws.Range("C22").Value2 = "8373999"
$ws.Range("D22").Value2 = "ELVER JOSE SEHUANES BULLOSO"
$ws.Range("E22").Value2 = "1911"
$ws.Range("F22").Value2 = 33125
$ws.Range("G22").Value2 = 828116
$ws.Range("B23").Value2 = "CC"
$ws.Range("C23").Value2 = "8373999"
$ws.Range("D23").Value2 = "ELVER JOSE SEHUANES BULLOSO"
$ws.Range("E23").Value2 = "1910"
$ws.Range("F23").Value2 = 33125
$ws.Range("G23").Value2 = 828116
$ws.Range("B24").Value2 = "CC"
$ws.Range("C24").Value2 = "8373999"
$ws.Range("D24").Value2 = "ELVER JOSE SEHUANES BULLOSO"
$ws.Range("E24").Value2 = "1909"
$ws.Range("F24").Value2 = 33125
$ws.Range("G24").Value2 = 828116
$ws.Range("B25").Value2 = "CC"
$ws.Range("C25").Value2 = "1051444079"
$ws.Range("D25").Value2 = "MILTON DAVID CERVANTES DE LA ROSA"
$ws.Range("E25").Value2 = "2003"
$ws.Range("F25").Value2 = 32021
$ws.Range("G25").Value2 = 828116
$ws.Range("B26").Value2 = "CC"
$ws.Range("C26").Value2 = "1051444079"
$ws.Range("D26").Value2 = "MILTON DAVID CERVANTES DE LA ROSA"
$ws.Range("E26").Value2 = "2002"
$ws.Range("F26").Value2 = 33125
$ws.Range("G26").Value2 = 828116
$ws.Range("B27").Value2 = "CC"
$ws.Range("C27").Value2 = "1051444079"
$ws.Range("D27").Value2 = "MILTON DAVID CERVANTES DE LA ROSA"
$ws.Range("E27").Value2 = "2001"
$ws.Range("F27").Value2 = 33125
$ws.Range("G27").Value2 = 828116
$ws.Range("B28").Value2 = "CC"
$ws.Range("C28").Value2 = "1051444079"
$ws.Range("D28").Value2 = "MILTON DAVID CERVANTES DE LA ROSA"
$ws.Range("E28").Value2 = "1912"
$ws.Range("F28").Value2 = 33125
$ws.Range("G28").Value2 = 828116
$ws.Range("B29").Value2 = "CC"
$ws.Range("C29").Value2 = "1051444079"
$ws.Range("D29").Value2 = "MILTON DAVID CERVANTES DE LA ROSA"
$ws.Range("E29").Value2 = "1911"
$ws.Range("F29").Value2 = 33125
$ws.Range("G29").Value2 = 828116
$ws.Range("B30").Value2 = "CC"
$ws.Range("C30").Value2 = "1051444079"
$ws.Range("D30").Value2 = "MILTON DAVID CERVANTES DE LA ROSA"
$ws.Range("E30").Value2 = "1910"
$ws.Range("F30").Value2 = 33125
$ws.Range("G30").Value2 = 828116
$ws.Range("B31").Value2 = "CC"
$ws.Range("C31").Value2 = "1051444079"
$ws.Range("D31").Value2 = "MILTON DAVID CERVANTES DE LA ROSA"
$ws.Range("E31").Value2 = "1909"
$ws.Range("F31").Value2 = 33125
$ws.Range("G31").Value2 = 828116
$ws.Range("B32").Value2 = "CC"
$ws.Range("C32").Value2 = "9146650"
$ws.Range("D32").Value2 = "LIBARDO ANTONIO CONTRERAS MARTINEZ"
$ws.Range("E32").Value2 = "2003"
$ws.Range("F32").Value2 = 32021
$ws.Range("G32").Value2 = 828116
$ws.Range("B33").Value2 = "CC"
$ws.Range("C33").Value2 = "9146650"
$ws.Range("D33").Value2 = "LIBARDO ANTONIO CONTRERAS MARTINEZ"
$ws.Range("E33").Value2 = "2002"
$ws.Range("F33").Value2 = 33125
$ws.Range("G33").Value2 = 828116
$ws.Range("B34").Value2 = "CC"
$ws.Range("C34").Value2 = "9146650"
$ws.Range("D34").Value2 = "LIBARDO ANTONIO CONTRERAS MARTINEZ"
$ws.Range("E34").Value2 = "2001"
$ws.Range("F34").Value2 = 33125
$ws.Range("G34").Value2 = 828116
$ws.Range("B35").Value2 = "CC"
$ws.Range("C35").Value2 = "9146650"
$ws.Range("D35").Value2 = "LIBARDO ANTONIO CONTRERAS MARTINEZ"
$ws.Range("E35").Value2 = "1910"
$ws.Range("F35").Value2 = 33125
$ws.Range("G35").Value2 = 828116
$ws.Range("B36").Value2 = "CC"
$ws.Range("C36").Value2 = "9146650"
$ws.Range("D36").Value2 = "LIBARDO ANTONIO CONTRERAS MARTINEZ"
$ws.Range("E36").Value2 = "1909"
$ws.Range("F36").Value2 = 33125
$ws.Range("G36").Value2 = 828116

# Step 5: Update the total "Valor Mora" shown at the top of the statement
$ws.Range("E11").Value2 = 692313
